$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: find the existing Hyperlink object attached to a given cell address
# (e.g. "$A$2") on a worksheet. Range-scoped .Hyperlinks / .Item() indexing
# are unreliable in this COM host, so walk the worksheet-level collection.
# ---------------------------------------------------------------------------
function Get-HyperlinkAt {
    param($Worksheet, [string]$CellAddress)

    foreach ($hl in $Worksheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $CellAddress) {
            return $hl
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Apply the "report generated for handback" update to one locale sheet.
# ---------------------------------------------------------------------------
function Update-LocaleSheet {
    param($SheetName, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status text: the report now reflects that the handback happened and the
    # target content is back in sync with en-US.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    # Capture the existing hyperlink URLs/display text for the source file
    # (A2) and the handoff target file (C2) before adding the new ones.
    $srcLink = Get-HyperlinkAt $ws "`$A`$2"
    $xlfLink = Get-HyperlinkAt $ws "`$C`$2"

    # New "Latest Target File" (E2) and "Latest Handback File" (F2) columns,
    # pointing at the same files that were handed off (A2 / C2), now marked
    # as the latest target/handback files. Give them the same look (underline,
    # cornflower-blue) as the other hyperlink cells on the row.
    $ws.Hyperlinks.Add($ws.Range("E2"), $srcLink.Address, "", "", $srcLink.TextToDisplay) | Out-Null
    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfLink.Address, "", "", $xlfLink.TextToDisplay) | Out-Null
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276

    # "Latest Handback DateTime" (G2) now records when the handback happened.
    $ws.Range("G2").Value = $HandbackDateTime
}

Update-LocaleSheet "zh-cn" "2016-01-26 09:34:35"
Update-LocaleSheet "de-de" "2016-01-26 09:34:56"
